# Updated cryptos list (Price / Volume(1h)) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.233.05"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.901.48"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D5").Value = "'306.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").Value = "'0.5337"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").Value = "'0.3812"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("D9").Value = "'0.07287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Value = "'22.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.06%  "

$ws.Range("D11").Value = "'0.9017"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'0.08207"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "'95.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "

$ws.Range("D14").Value = "'5.334"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").Value = "'14.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").Value = "27.263.13"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").Value = "'5.028"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "1.085.30"
$ws.Range("E21").Value = "  -42.88%  "

$ws.Range("D22").Value = "'10.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").Value = "'6.509"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "

$ws.Range("D24").Value = "'149.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.60%  "

$ws.Range("D25").Value = "'2.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").Value = "'18.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "'1.747"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").Value = "'116.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").Value = "'4.811"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").Value = "'0.09241"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").Value = "'0.8304"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'0.05056"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  -1.36%  "

$ws.Range("D35").Value = "'3.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.30%  "

$ws.Range("D36").Value = "'3.334"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.20%  "

$ws.Range("D37").Value = "'2.676"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "

$ws.Range("D38").Value = "'0.5742"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").Value = "'0.02004"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("D40").Value = "'1.074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").Value = "'9.358"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("D43").Value = "'117.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "

$ws.Range("D44").Value = "'0.1521"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").Value = "'0.4940"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'10.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "

$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("D49").Value = "'38.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "

$ws.Range("D50").Value = "'0.06147"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("D51").Value = "'63.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
